{"js": "const pairs = [\n  [\"2024-10-04 Friday\", \"2024-10-05 Saturday\"],\n  [\"948\u00f76=158, 0\", \"769\u00f72=384, 1\"],\n  [\"673\u00f76=112, 1\", \"132\u00f76=22, 0\"],\n  [\"125\u00f73=41, 2\", \"836\u00f79=92, 8\"],\n  [\"883\u00f76=147, 1\", \"789\u00f75=157, 4\"],\n  [\"841\u00f74=210, 1\", \"281\u00f77=40, 1\"],\n  [\"276\u00f75=55, 1\", \"895\u00f75=179, 0\"],\n  [\"961\u00f79=106, 7\", \"356\u00f79=39, 5\"],\n  [\"136\u00f77=19, 3\", \"885\u00f73=295, 0\"],\n  [\"621\u00f74=155, 1\", \"961\u00f75=192, 1\"],\n  [\"813\u00f77=116, 1\", \"877\u00f78=109, 5\"],\n  [\"244\u00f75=48, 4\", \"334\u00f74=83, 2\"],\n  [\"435\u00f77=62, 1\", \"336\u00f76=56, 0\"],\n  [\"566\u00f73=188, 2\", \"585\u00f76=97, 3\"],\n  [\"688\u00f78=86, 0\", \"540\u00f72=270, 0\"],\n  [\"521\u00f76=86, 5\", \"949\u00f79=105, 4\"],\n  [\"765\u00f74=191, 1\", \"306\u00f73=102, 0\"],\n  [\"513\u00f74=128, 1\", \"415\u00f73=138, 1\"],\n  [\"143\u00f73=47, 2\", \"642\u00f72=321, 0\"],\n  [\"910\u00f78=113, 6\", \"238\u00f75=47, 3\"],\n  [\"919\u00f77=131, 2\", \"676\u00f75=135, 1\"],\n  [\"108\u00f78=13, 4\", \"621\u00f73=207, 0\"],\n  [\"751\u00f76=125, 1\", \"539\u00f74=134, 3\"],\n  [\"235\u00f73=78, 1\", \"793\u00f73=264, 1\"],\n  [\"981\u00f72=490, 1\", \"975\u00f75=195, 0\"],\n  [\"602\u00f77=86, 0\", \"881\u00f77=125, 6\"],\n];\n\nfor (const [oldText, newText] of pairs) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2024-10-04 Friday\", \"2024-10-05 Saturday\"),\n    @(\"948\u00f76=158, 0\", \"769\u00f72=384, 1\"),\n    @(\"673\u00f76=112, 1\", \"132\u00f76=22, 0\"),\n    @(\"125\u00f73=41, 2\", \"836\u00f79=92, 8\"),\n    @(\"883\u00f76=147, 1\", \"789\u00f75=157, 4\"),\n    @(\"841\u00f74=210, 1\", \"281\u00f77=40, 1\"),\n    @(\"276\u00f75=55, 1\", \"895\u00f75=179, 0\"),\n    @(\"961\u00f79=106, 7\", \"356\u00f79=39, 5\"),\n    @(\"136\u00f77=19, 3\", \"885\u00f73=295, 0\"),\n    @(\"621\u00f74=155, 1\", \"961\u00f75=192, 1\"),\n    @(\"813\u00f77=116, 1\", \"877\u00f78=109, 5\"),\n    @(\"244\u00f75=48, 4\", \"334\u00f74=83, 2\"),\n    @(\"435\u00f77=62, 1\", \"336\u00f76=56, 0\"),\n    @(\"566\u00f73=188, 2\", \"585\u00f76=97, 3\"),\n    @(\"688\u00f78=86, 0\", \"540\u00f72=270, 0\"),\n    @(\"521\u00f76=86, 5\", \"949\u00f79=105, 4\"),\n    @(\"765\u00f74=191, 1\", \"306\u00f73=102, 0\"),\n    @(\"513\u00f74=128, 1\", \"415\u00f73=138, 1\"),\n    @(\"143\u00f73=47, 2\", \"642\u00f72=321, 0\"),\n    @(\"910\u00f78=113, 6\", \"238\u00f75=47, 3\"),\n    @(\"919\u00f77=131, 2\", \"676\u00f75=135, 1\"),\n    @(\"108\u00f78=13, 4\", \"621\u00f73=207, 0\"),\n    @(\"751\u00f76=125, 1\", \"539\u00f74=134, 3\"),\n    @(\"235\u00f73=78, 1\", \"793\u00f73=264, 1\"),\n    @(\"981\u00f72=490, 1\", \"975\u00f75=195, 0\"),\n    @(\"602\u00f77=86, 0\", \"881\u00f77=125, 6\"),\n)\n\nforeach ($pair in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute($pair[0], $false, $false, $false, $false, $false, $true, 1, $false, $pair[1], 2)\n}"}
